$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.158102766798419
$ws.Range("C2").Value = 0.6363636363636364
$ws.Range("J2").Value = 0.0158102766798419
$ws.Range("P2").Value = 0.1185770750988142
$ws.Range("S2").Value = 0.07114624505928854
$ws.Range("C3").Value = 0.05813953488372093
$ws.Range("J3").Value = 0.01744186046511628
$ws.Range("P3").Value = 0.7732558139534884
$ws.Range("S3").Value = 0.1511627906976744
$ws.Range("O4").Value = 0.02272727272727273
$ws.Range("P4").Value = 0.7727272727272727
$ws.Range("S4").Value = 0.2045454545454546
$ws.Range("B6").Value = 0.07555555555555556
$ws.Range("D6").Value = 0.01777777777777778
$ws.Range("F6").Value = 0.09333333333333334
$ws.Range("J6").Value = 0.2533333333333334
$ws.Range("O6").Value = 0.01777777777777778
$ws.Range("Q6").Value = 0.09777777777777778
$ws.Range("R6").Value = 0.06222222222222222
$ws.Range("S6").Value = 0.3822222222222222
$ws.Range("B7").Value = 0.08641975308641975
$ws.Range("D7").Value = 0.02469135802469136
$ws.Range("F7").Value = 0.04938271604938271
$ws.Range("J7").Value = 0.1049382716049383
$ws.Range("O7").Value = 0.01851851851851852
$ws.Range("Q7").Value = 0.228395061728395
$ws.Range("R7").Value = 0.09876543209876543
$ws.Range("B8").Value = 0.08531746031746032
$ws.Range("D8").Value = 0.01388888888888889
$ws.Range("E8").Value = 0.001984126984126984
$ws.Range("F8").Value = 0.05158730158730158
$ws.Range("J8").Value = 0.08928571428571429
$ws.Range("O8").Value = 0.01587301587301587
$ws.Range("Q8").Value = 0.2063492063492063
$ws.Range("R8").Value = 0.1071428571428571
$ws.Range("S8").Value = 0.4285714285714285
$ws.Range("B9").Value = 0.08056872037914692
$ws.Range("D9").Value = 0.03317535545023697
$ws.Range("F9").Value = 0.05213270142180094
$ws.Range("J9").Value = 0.08530805687203792
$ws.Range("O9").Value = 0.02369668246445497
$ws.Range("Q9").Value = 0.1943127962085308
$ws.Range("R9").Value = 0.1611374407582938
$ws.Range("S9").Value = 0.3696682464454976
$ws.Range("B10").Value = 0.09041309431021044
$ws.Range("D10").Value = 0.01792673421667966
$ws.Range("E10").Value = 0.000779423226812159
$ws.Range("F10").Value = 0.06780982073265783
$ws.Range("J10").Value = 0.1239282930631333
$ws.Range("O10").Value = 0.01480904130943102
$ws.Range("Q10").Value = 0.2213561964146532
$ws.Range("R10").Value = 0.08729540140296181
$ws.Range("S10").Value = 0.3756819953234606
$ws.Range("G11").Value = 0.1347826086956522
$ws.Range("J11").Value = 0.1130434782608696
$ws.Range("K11").Value = 0.2217391304347826
$ws.Range("L11").Value = 0.5260869565217391
$ws.Range("S11").Value = 0.004347826086956522
$ws.Range("G12").Value = 0.7983870967741935
$ws.Range("J12").Value = 0.1290322580645161
$ws.Range("L12").Value = 0.04838709677419355
$ws.Range("S12").Value = 0.02419354838709677
$ws.Range("G13").Value = 0.7755102040816326
$ws.Range("J13").Value = 0.2244897959183673
$ws.Range("F15").Value = 0.02347417840375587
$ws.Range("H15").Value = 0.1643192488262911
$ws.Range("I15").Value = 0.04694835680751173
$ws.Range("J15").Value = 0.3849765258215962
$ws.Range("K15").Value = 0.07511737089201878
$ws.Range("M15").Value = 0.009389671361502348
$ws.Range("O15").Value = 0.0892018779342723
$ws.Range("S15").Value = 0.2065727699530517
$ws.Range("F16").Value = 0.005235602094240838
$ws.Range("H16").Value = 0.2094240837696335
$ws.Range("I16").Value = 0.06806282722513089
$ws.Range("J16").Value = 0.4240837696335079
$ws.Range("K16").Value = 0.08376963350785341
$ws.Range("M16").Value = 0.02617801047120419
$ws.Range("O16").Value = 0.04712041884816754
$ws.Range("S16").Value = 0.1361256544502618
$ws.Range("F17").Value = 0.02066115702479339
$ws.Range("H17").Value = 0.1632231404958678
$ws.Range("I17").Value = 0.115702479338843
$ws.Range("J17").Value = 0.4566115702479339
$ws.Range("K17").Value = 0.05165289256198347
$ws.Range("M17").Value = 0.02272727272727273
$ws.Range("O17").Value = 0.0640495867768595
$ws.Range("S17").Value = 0.1053719008264463
$ws.Range("F18").Value = 0.02586206896551724
$ws.Range("H18").Value = 0.2112068965517241
$ws.Range("I18").Value = 0.103448275862069
$ws.Range("J18").Value = 0.4094827586206897
$ws.Range("K18").Value = 0.06896551724137931
$ws.Range("M18").Value = 0.03448275862068965
$ws.Range("O18").Value = 0.04741379310344827
$ws.Range("S18").Value = 0.09913793103448276
$ws.Range("F19").Value = 0.01885310290652003
$ws.Range("H19").Value = 0.2388059701492537
$ws.Range("I19").Value = 0.08719560094265515
$ws.Range("J19").Value = 0.3637077769049489
$ws.Range("K19").Value = 0.08248232521602514
$ws.Range("M19").Value = 0.01885310290652003
$ws.Range("O19").Value = 0.06284367635506677
$ws.Range("S19").Value = 0.1272584446190102
